{"js": "// Update the date line and the 99 changed arithmetic-problem cells in the\n// 20x5 table. One cell (\"81-37=\" at row 7 / col 5, 0-indexed row 6 col 4)\n// is left untouched, matching the source diff.\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\n// --- 1) Update the date heading (first paragraph of the document) ---\nconst firstPara = body.paragraphs.items[0];\nfirstPara.load(\"text\");\nawait context.sync();\nif (firstPara.text.trim() === \"2024-12-05 Thursday\") {\n  firstPara.insertText(\"2024-12-06 Friday\", \"Replace\");\n}\n\n// --- 2) Update the arithmetic table cells ---\nconst replacements = [\n  [0, 0, \"67+25=\", \"70-45=\"],\n  [0, 1, \"87-9=\", \"84-5=\"],\n  [0, 2, \"24+38=\", \"19+75=\"],\n  [0, 3, \"75-58=\", \"39+18=\"],\n  [0, 4, \"27-9=\", \"97-79=\"],\n  [1, 0, \"59+9=\", \"55+26=\"],\n  [1, 1, \"68+17=\", \"58+14=\"],\n  [1, 2, \"26+29=\", \"87-69=\"],\n  [1, 3, \"34+38=\", \"86+6=\"],\n  [1, 4, \"44+27=\", \"18+36=\"],\n  [2, 0, \"25+9=\", \"64-25=\"],\n  [2, 1, \"71-57=\", \"68+8=\"],\n  [2, 2, \"94-57=\", \"84+9=\"],\n  [2, 3, \"61-9=\", \"19+7=\"],\n  [2, 4, \"41-7=\", \"50-7=\"],\n  [3, 0, \"51-18=\", \"78+17=\"],\n  [3, 1, \"6+75=\", \"27+45=\"],\n  [3, 2, \"41-35=\", \"15+6=\"],\n  [3, 3, \"66+26=\", \"32+9=\"],\n  [3, 4, \"8+56=\", \"26+15=\"],\n  [4, 0, \"19+79=\", \"77+17=\"],\n  [4, 1, \"58+5=\", \"92-33=\"],\n  [4, 2, \"52-23=\", \"57+9=\"],\n  [4, 3, \"6+29=\", \"32+39=\"],\n  [4, 4, \"91-37=\", \"48+33=\"],\n  [5, 0, \"51-39=\", \"6+7=\"],\n  [5, 1, \"52-28=\", \"51-14=\"],\n  [5, 2, \"25+28=\", \"39+32=\"],\n  [5, 3, \"57+15=\", \"91-6=\"],\n  [5, 4, \"9+3=\", \"72-24=\"],\n  [6, 0, \"19+12=\", \"28+39=\"],\n  [6, 1, \"24+69=\", \"34-28=\"],\n  [6, 2, \"65-56=\", \"25+39=\"],\n  [6, 3, \"8+29=\", \"27+34=\"],\n  [7, 0, \"28+34=\", \"9+48=\"],\n  [7, 1, \"94-38=\", \"26+49=\"],\n  [7, 2, \"80-49=\", \"60-33=\"],\n  [7, 3, \"13+59=\", \"84-15=\"],\n  [7, 4, \"8+13=\", \"19+78=\"],\n  [8, 0, \"19+28=\", \"14+17=\"],\n  [8, 1, \"36+26=\", \"88-79=\"],\n  [8, 2, \"83-24=\", \"36+15=\"],\n  [8, 3, \"90-9=\", \"8+45=\"],\n  [8, 4, \"85-48=\", \"81-13=\"],\n  [9, 0, \"6+38=\", \"49+3=\"],\n  [9, 1, \"44+48=\", \"17+38=\"],\n  [9, 2, \"18+65=\", \"70-59=\"],\n  [9, 3, \"46-37=\", \"22-17=\"],\n  [9, 4, \"83-79=\", \"7+54=\"],\n  [10, 0, \"85-47=\", \"54+37=\"],\n  [10, 1, \"9+77=\", \"92-15=\"],\n  [10, 2, \"91-29=\", \"59+23=\"],\n  [10, 3, \"74-55=\", \"92-87=\"],\n  [10, 4, \"4+17=\", \"90-35=\"],\n  [11, 0, \"79+15=\", \"9+4=\"],\n  [11, 1, \"8+88=\", \"77-29=\"],\n  [11, 2, \"70-16=\", \"38+49=\"],\n  [11, 3, \"35-8=\", \"9+53=\"],\n  [11, 4, \"69+5=\", \"48+9=\"],\n  [12, 0, \"37+18=\", \"64+9=\"],\n  [12, 1, \"5+28=\", \"19+63=\"],\n  [12, 2, \"36+27=\", \"80-1=\"],\n  [12, 3, \"85-7=\", \"52-45=\"],\n  [12, 4, \"88+3=\", \"54-18=\"],\n  [13, 0, \"57-9=\", \"17+64=\"],\n  [13, 1, \"9+42=\", \"80-23=\"],\n  [13, 2, \"5+57=\", \"91-37=\"],\n  [13, 3, \"19+7=\", \"51-49=\"],\n  [13, 4, \"36+55=\", \"29+12=\"],\n  [14, 0, \"8+13=\", \"5+77=\"],\n  [14, 1, \"16+79=\", \"19+44=\"],\n  [14, 2, \"9+67=\", \"44-15=\"],\n  [14, 3, \"90-16=\", \"28+54=\"],\n  [14, 4, \"19+53=\", \"27-8=\"],\n  [15, 0, \"17+69=\", \"83-54=\"],\n  [15, 1, \"85-6=\", \"12-4=\"],\n  [15, 2, \"73-19=\", \"70-46=\"],\n  [15, 3, \"51-6=\", \"13+58=\"],\n  [15, 4, \"29+37=\", \"35-18=\"],\n  [16, 0, \"40-9=\", \"47+44=\"],\n  [16, 1, \"18+59=\", \"50-33=\"],\n  [16, 2, \"7+7=\", \"67+25=\"],\n  [16, 3, \"68+3=\", \"8+19=\"],\n  [16, 4, \"71-32=\", \"48+28=\"],\n  [17, 0, \"25+19=\", \"55-27=\"],\n  [17, 1, \"25+49=\", \"16+38=\"],\n  [17, 2, \"82-4=\", \"8+26=\"],\n  [17, 3, \"5+27=\", \"27+35=\"],\n  [17, 4, \"65-29=\", \"37-29=\"],\n  [18, 0, \"51-28=\", \"92-67=\"],\n  [18, 1, \"46+16=\", \"42-7=\"],\n  [18, 2, \"11-2=\", \"34+8=\"],\n  [18, 3, \"53+38=\", \"4+58=\"],\n  [18, 4, \"67+6=\", \"70-1=\"],\n  [19, 0, \"89+6=\", \"63-57=\"],\n  [19, 1, \"55+28=\", \"71-56=\"],\n  [19, 2, \"83+9=\", \"75+8=\"],\n  [19, 3, \"26+59=\", \"70-34=\"],\n  [19, 4, \"18+76=\", \"78-19=\"],\n];\n\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nfor (const [row, col, oldText, newText] of replacements) {\n  const cell = table.getCell(row, col);\n  cell.body.paragraphs.load(\"items\");\n  await context.sync();\n  const para = cell.body.paragraphs.items[0];\n  para.load(\"text\");\n  await context.sync();\n  if (para.text.trim() === oldText) {\n    para.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 99 changed arithmetic-problem cells in the\n# 20x5 table. One cell (\"81-37=\" at Row 7 / Col 5) is left untouched,\n# matching the source diff.\n\n$d = $word.ActiveDocument\n\n# --- 1) Update the date heading (first paragraph of the document) ---\n$firstPara = $d.Paragraphs.Item(1)\nif ($firstPara.Range.Text.TrimEnd([char]13, [char]7) -eq \"2024-12-05 Thursday\") {\n    $firstPara.Range.Text = \"2024-12-06 Friday\"\n}\n\n# --- 2) Update the arithmetic table cells ---\n$replacements = @(\n    @{Row=1; Col=1; Old='67+25='; New='70-45='},\n    @{Row=1; Col=2; Old='87-9='; New='84-5='},\n    @{Row=1; Col=3; Old='24+38='; New='19+75='},\n    @{Row=1; Col=4; Old='75-58='; New='39+18='},\n    @{Row=1; Col=5; Old='27-9='; New='97-79='},\n    @{Row=2; Col=1; Old='59+9='; New='55+26='},\n    @{Row=2; Col=2; Old='68+17='; New='58+14='},\n    @{Row=2; Col=3; Old='26+29='; New='87-69='},\n    @{Row=2; Col=4; Old='34+38='; New='86+6='},\n    @{Row=2; Col=5; Old='44+27='; New='18+36='},\n    @{Row=3; Col=1; Old='25+9='; New='64-25='},\n    @{Row=3; Col=2; Old='71-57='; New='68+8='},\n    @{Row=3; Col=3; Old='94-57='; New='84+9='},\n    @{Row=3; Col=4; Old='61-9='; New='19+7='},\n    @{Row=3; Col=5; Old='41-7='; New='50-7='},\n    @{Row=4; Col=1; Old='51-18='; New='78+17='},\n    @{Row=4; Col=2; Old='6+75='; New='27+45='},\n    @{Row=4; Col=3; Old='41-35='; New='15+6='},\n    @{Row=4; Col=4; Old='66+26='; New='32+9='},\n    @{Row=4; Col=5; Old='8+56='; New='26+15='},\n    @{Row=5; Col=1; Old='19+79='; New='77+17='},\n    @{Row=5; Col=2; Old='58+5='; New='92-33='},\n    @{Row=5; Col=3; Old='52-23='; New='57+9='},\n    @{Row=5; Col=4; Old='6+29='; New='32+39='},\n    @{Row=5; Col=5; Old='91-37='; New='48+33='},\n    @{Row=6; Col=1; Old='51-39='; New='6+7='},\n    @{Row=6; Col=2; Old='52-28='; New='51-14='},\n    @{Row=6; Col=3; Old='25+28='; New='39+32='},\n    @{Row=6; Col=4; Old='57+15='; New='91-6='},\n    @{Row=6; Col=5; Old='9+3='; New='72-24='},\n    @{Row=7; Col=1; Old='19+12='; New='28+39='},\n    @{Row=7; Col=2; Old='24+69='; New='34-28='},\n    @{Row=7; Col=3; Old='65-56='; New='25+39='},\n    @{Row=7; Col=4; Old='8+29='; New='27+34='},\n    @{Row=8; Col=1; Old='28+34='; New='9+48='},\n    @{Row=8; Col=2; Old='94-38='; New='26+49='},\n    @{Row=8; Col=3; Old='80-49='; New='60-33='},\n    @{Row=8; Col=4; Old='13+59='; New='84-15='},\n    @{Row=8; Col=5; Old='8+13='; New='19+78='},\n    @{Row=9; Col=1; Old='19+28='; New='14+17='},\n    @{Row=9; Col=2; Old='36+26='; New='88-79='},\n    @{Row=9; Col=3; Old='83-24='; New='36+15='},\n    @{Row=9; Col=4; Old='90-9='; New='8+45='},\n    @{Row=9; Col=5; Old='85-48='; New='81-13='},\n    @{Row=10; Col=1; Old='6+38='; New='49+3='},\n    @{Row=10; Col=2; Old='44+48='; New='17+38='},\n    @{Row=10; Col=3; Old='18+65='; New='70-59='},\n    @{Row=10; Col=4; Old='46-37='; New='22-17='},\n    @{Row=10; Col=5; Old='83-79='; New='7+54='},\n    @{Row=11; Col=1; Old='85-47='; New='54+37='},\n    @{Row=11; Col=2; Old='9+77='; New='92-15='},\n    @{Row=11; Col=3; Old='91-29='; New='59+23='},\n    @{Row=11; Col=4; Old='74-55='; New='92-87='},\n    @{Row=11; Col=5; Old='4+17='; New='90-35='},\n    @{Row=12; Col=1; Old='79+15='; New='9+4='},\n    @{Row=12; Col=2; Old='8+88='; New='77-29='},\n    @{Row=12; Col=3; Old='70-16='; New='38+49='},\n    @{Row=12; Col=4; Old='35-8='; New='9+53='},\n    @{Row=12; Col=5; Old='69+5='; New='48+9='},\n    @{Row=13; Col=1; Old='37+18='; New='64+9='},\n    @{Row=13; Col=2; Old='5+28='; New='19+63='},\n    @{Row=13; Col=3; Old='36+27='; New='80-1='},\n    @{Row=13; Col=4; Old='85-7='; New='52-45='},\n    @{Row=13; Col=5; Old='88+3='; New='54-18='},\n    @{Row=14; Col=1; Old='57-9='; New='17+64='},\n    @{Row=14; Col=2; Old='9+42='; New='80-23='},\n    @{Row=14; Col=3; Old='5+57='; New='91-37='},\n    @{Row=14; Col=4; Old='19+7='; New='51-49='},\n    @{Row=14; Col=5; Old='36+55='; New='29+12='},\n    @{Row=15; Col=1; Old='8+13='; New='5+77='},\n    @{Row=15; Col=2; Old='16+79='; New='19+44='},\n    @{Row=15; Col=3; Old='9+67='; New='44-15='},\n    @{Row=15; Col=4; Old='90-16='; New='28+54='},\n    @{Row=15; Col=5; Old='19+53='; New='27-8='},\n    @{Row=16; Col=1; Old='17+69='; New='83-54='},\n    @{Row=16; Col=2; Old='85-6='; New='12-4='},\n    @{Row=16; Col=3; Old='73-19='; New='70-46='},\n    @{Row=16; Col=4; Old='51-6='; New='13+58='},\n    @{Row=16; Col=5; Old='29+37='; New='35-18='},\n    @{Row=17; Col=1; Old='40-9='; New='47+44='},\n    @{Row=17; Col=2; Old='18+59='; New='50-33='},\n    @{Row=17; Col=3; Old='7+7='; New='67+25='},\n    @{Row=17; Col=4; Old='68+3='; New='8+19='},\n    @{Row=17; Col=5; Old='71-32='; New='48+28='},\n    @{Row=18; Col=1; Old='25+19='; New='55-27='},\n    @{Row=18; Col=2; Old='25+49='; New='16+38='},\n    @{Row=18; Col=3; Old='82-4='; New='8+26='},\n    @{Row=18; Col=4; Old='5+27='; New='27+35='},\n    @{Row=18; Col=5; Old='65-29='; New='37-29='},\n    @{Row=19; Col=1; Old='51-28='; New='92-67='},\n    @{Row=19; Col=2; Old='46+16='; New='42-7='},\n    @{Row=19; Col=3; Old='11-2='; New='34+8='},\n    @{Row=19; Col=4; Old='53+38='; New='4+58='},\n    @{Row=19; Col=5; Old='67+6='; New='70-1='},\n    @{Row=20; Col=1; Old='89+6='; New='63-57='},\n    @{Row=20; Col=2; Old='55+28='; New='71-56='},\n    @{Row=20; Col=3; Old='83+9='; New='75+8='},\n    @{Row=20; Col=4; Old='26+59='; New='70-34='},\n    @{Row=20; Col=5; Old='18+76='; New='78-19='}\n)\n\n$table = $d.Tables.Item(1)\n\nforeach ($r in $replacements) {\n    $cell = $table.Cell($r.Row, $r.Col)\n    $cellText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($cellText -eq $r.Old) {\n        $cell.Range.Text = $r.New\n    }\n}\n\nWrite-Output \"done\"\n"}
